$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Capture the comment text that currently lives on H1/I1 before the columns
# shift, so the surviving "Tag" / "Instrument *" columns keep their notes.
$tagComment = $ws.Range("H1").Comment.Text()
$instrumentComment = $ws.Range("I1").Comment.Text()

# Remove the "Type *" (F) and "Folio No" (G) columns entirely - this shifts
# Tag/Instrument left from H/I to F/G along with the rest of the data.
$ws.Range("F:G").EntireColumn.Delete()

# The engine doesn't relocate cell comments when columns are deleted, so
# refresh F1/G1's comment text to match what previously lived on H1/I1, and
# drop the now-orphaned comments left behind past the new right edge.
$ws.Range("F1").ClearComments()
$ws.Range("F1").AddComment($tagComment)

$ws.Range("G1").ClearComments()
$ws.Range("G1").AddComment($instrumentComment)

$ws.Range("H1").Comment.Delete()
$ws.Range("I1").Comment.Delete()

# Match the post-edit selection state (whole column F selected).
[void]$ws.Range("F:F").Select()
